# Commit: Wed, Apr 29, 2020  2:06:49 PM
#
# The table on slide 16 (the "Total Outflow / ..." summary table placed
# below the screen-clipping picture) had a different table style applied
# to it - switch it from the deck's custom "Table_0" style to the
# built-in PowerPoint table style {0E074CAA-A7D7-42B9-8F07-CF459D60DF05}.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(16)

# Shape order on this slide: 1 = title placeholder, 2 = screen-clipping
# picture, 3 = the graphicFrame holding the table.
$sh = $s.Shapes.Item(3)

$tbl = $sh.Table
$tbl.ApplyStyle("{0E074CAA-A7D7-42B9-8F07-CF459D60DF05}")
